# Aggiunta la citta' ai nomi degli ospedali
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cent")

$ws.Range("B2").Value = "Ospedale Santa Maria (Bari)"
$ws.Range("B3").Value = "Ospedale San Paolo (Bari)"
$ws.Range("B4").Value = "Ospedale Papa Giovanni XXXIII (Bari)"
$ws.Range("B5").Value = "Ospedale Don Tonino Bello (Molfetta)"
$ws.Range("B6").Value = "Ospedale Monsignor Raffaele di Miccoli (Barletta)"
$ws.Range("B7").Value = "Ospedale IRCSS Saverio De Bellis (Castellana)"

$ws.Range("B7").Select() | Out-Null
